$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing cells flagged by the new structure-detection pass ---
$ws.Range("Q62").Value = 0
$ws.Range("Q68").Value = 0
$ws.Range("Q79").Value = 0
$ws.Range("O399").Value = 1
$ws.Range("R399").Value = 0
$ws.Range("R400").Value = 0

# --- Append newly scraped weekly rows 401-432 ---
# Columns: Datetime Open High Low Close AdjClose Volume Year Month Day Hour Minute Second Week isPivot two_line_structure detect_structure backup
# "NA" marks columns left blank in the source (Adj Close / backup).
$data = @"
401 45474 651.452870649912 653.0515036076249 629.4713320697002 645.1181640625 NA 13526077 2024 7 1 0 0 0 27 0 0 0 NA
402 45481 645.4179573923491 662.5435278217709 625.5346323064418 633.2681274414062 NA 17778019 2024 7 8 0 0 0 28 0 0 0 NA
403 45488 631.4696261366805 657.327912541372 615.0833646091797 623.3564453125 NA 18414538 2024 7 15 0 0 0 29 0 0 0 NA
404 45495 623.4763509919328 676.6317187540581 602.593865483404 670.09716796875 NA 26492118 2024 7 22 0 0 0 30 0 0 0 NA
405 45502 678.6300658439475 680.5484376126768 620.4588890743815 628.4122314453125 NA 53721686 2024 7 29 0 0 0 31 0 0 2 NA
406 45509 615.2832065641397 619.2798497519341 581.5115838240654 602.93359375 NA 25161164 2024 8 5 0 0 0 32 0 0 0 NA
407 45516 598 601.1799926757812 565 568.739990234375 NA 25170530 2024 8 12 0 0 0 33 0 0 0 NA
408 45523 571.5999755859375 643.780029296875 567.1799926757812 632.8400268554688 NA 33732464 2024 8 19 0 0 0 34 0 0 0 NA
409 45530 630 635.9400024414062 597.5999755859375 600.3599853515625 NA 20911644 2024 8 26 0 0 0 35 0 0 0 NA
410 45537 600.3400268554688 617.5999755859375 591.6799926757812 597.2999877929688 NA 23815358 2024 9 2 0 0 0 36 0 0 0 NA
411 45544 600 668 549.219970703125 646.6500244140625 NA 34955823 2024 9 9 0 0 0 37 1 0 0 NA
412 45551 648 661.4500122070312 616.25 654.4500122070312 NA 32878196 2024 9 16 0 0 0 38 0 0 0 NA
413 45558 654.4500122070312 663 606 608.5499877929688 NA 21068665 2024 9 23 0 0 0 39 0 0 0 NA
414 45565 624 624 570.5 578.9000244140625 NA 30400469 2024 9 30 0 0 0 40 0 0 0 NA
415 45572 574.9000244140625 606.2000122070312 536.8499755859375 592.3499755859375 NA 49963574 2024 10 7 0 0 0 41 2 0 0 NA
416 45579 596.7999877929688 613.9000244140625 578.7000122070312 588.75 NA 20911424 2024 10 14 0 0 0 42 0 0 0 NA
417 45586 597.5999755859375 631.2000122070312 566.7000122070312 611.2000122070312 NA 40918789 2024 10 21 0 0 0 43 0 0 2 NA
418 45593 608.25 615.5499877929688 589.4000244140625 608.5 NA 16856270 2024 10 28 0 0 0 44 0 0 0 NA
419 45600 611 615.9500122070312 580.4500122070312 591.5499877929688 NA 24182068 2024 11 4 0 0 0 45 0 0 0 NA
420 45607 588 607.9500122070312 565 575.6500244140625 NA 18881902 2024 11 11 0 0 0 46 0 0 0 NA
421 45614 593 639.5999755859375 582.8499755859375 616.4000244140625 NA 33360939 2024 11 18 0 0 0 47 0 0 1 NA
422 45621 622 638.5999755859375 595 621.2000122070312 NA 30080570 2024 11 25 0 0 0 48 0 0 0 NA
423 45628 625 645 598.7999877929688 644.0499877929688 NA 40853205 2024 12 2 0 0 0 49 0 0 0 NA
424 45635 644 653 635.5 645.6500244140625 NA 23919780 2024 12 9 0 0 0 50 0 0 0 NA
425 45642 648.9500122070312 657.2000122070312 609.1500244140625 612.5499877929688 NA 27038702 2024 12 16 0 0 0 51 0 0 0 NA
426 45649 621.2999877929688 629.9000244140625 620 624.2999877929688 NA 9437354 2024 12 23 0 0 0 52 0 0 0 NA
427 45656 621.1500244140625 663.5999755859375 621 652.2000122070312 NA 32692157 2024 12 30 0 0 0 1 1 0 0 NA
428 45663 645 645.9500122070312 590.4000244140625 596.5499877929688 NA 20224638 2025 1 6 0 0 0 2 0 0 0 NA
429 45670 585 591.3499755859375 545.5499877929688 556.3499755859375 NA 31551033 2025 1 13 0 0 0 3 0 0 0 NA
430 45677 557 563 520 541.7000122070312 NA 28619486 2025 1 20 0 0 0 4 0 0 0 NA
431 45684 538.9000244140625 582.9000244140625 516.9500122070312 561.4000244140625 NA 24700124 2025 1 27 0 0 0 5 0 0 0 NA
432 45691 560 592.9500122070312 548.5499877929688 553.7000122070312 NA 22596034 2025 2 3 0 0 0 6 0 0 0 NA
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split "\s+"
    $rowNum = [int]$parts[0]
    $ws.Cells.Item($rowNum, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    for ($c = 0; $c -lt 17; $c++) {
        $token = $parts[$c + 1]
        if ($token -ne "NA") {
            $ws.Cells.Item($rowNum, $c + 1).Value = [double]$token
        }
    }
}

Write-Host "done"
